# Update "想去人数" (want-to-go count) values for two sheets: 展览 and 全部类型
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition): rows 2,4,5,6 in column F
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1090
$ws1.Range("F4").Value = 1682
$ws1.Range("F5").Value = 759
$ws1.Range("F6").Value = 185

# Sheet "全部类型" (All types): rows 2,4,6,7 in column F
$ws2 = $wb.Worksheets.Item("全部类型")
$ws2.Range("F2").Value = 1090
$ws2.Range("F4").Value = 1682
$ws2.Range("F6").Value = 759
$ws2.Range("F7").Value = 185
